$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.685.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.495.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.15"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.492.77"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.191"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.25"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.580"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.06"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000273"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.067.52"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "607.11"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.504.56"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.770.88"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.870"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -18.54%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.55"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.17"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "680.25"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +18.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.96"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.89"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0997"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.54"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.68"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0473"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.17%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.40"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.141"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.311.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.311"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.91"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.18"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0687"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.129"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.88"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.06%  "
